# Apply Price/Volume(1h) updates from the coinranking.com scrape refresh.
# Cells whose new text would otherwise be auto-parsed as a number by Excel
# (General-format column) are entered with a leading apostrophe so they
# stay plain text, matching the original inlineStr cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.830.30"
$ws.Range("D3").Value = "1.940.71"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'242.85"
$ws.Range("E5").Value = "  -1.11%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "'0.4892"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'0.2941"
$ws.Range("E8").Value = "  -0.54%  "
$ws.Range("D9").Value = "'0.06921"
$ws.Range("E9").Value = "  +1.47%  "
$ws.Range("D10").Value = "'19.37"
$ws.Range("E10").Value = "  +1.00%  "
$ws.Range("D11").Value = "'105.50"
$ws.Range("E11").Value = "  -1.00%  "
$ws.Range("D12").Value = "1.934.44"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").Value = "'0.07723"
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("D14").Value = "'5.367"
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("D15").Value = "'0.6997"
$ws.Range("E15").Value = "  -1.71%  "
$ws.Range("D16").Value = "'273.02"
$ws.Range("E16").Value = "  -4.46%  "
$ws.Range("D17").Value = "30.833.98"
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("D18").Value = "'0.000007722"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("E19").Value = "  -0.62%  "
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.195.75"
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "'5.522"
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "'6.556"
$ws.Range("E24").Value = "  -0.67%  "
$ws.Range("E25").Value = "  -1.77%  "
$ws.Range("D26").Value = "'166.51"
$ws.Range("E26").Value = "  -1.36%  "
$ws.Range("D27").Value = "'19.62"
$ws.Range("E27").Value = "  -1.75%  "
$ws.Range("D28").Value = "'2.173"
$ws.Range("E28").Value = "  -0.90%  "
$ws.Range("D29").Value = "'0.1039"
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("D30").Value = "'1.389"
$ws.Range("E30").Value = "  -3.41%  "
$ws.Range("D31").Value = "'4.583"
$ws.Range("E31").Value = "  -3.34%  "
$ws.Range("E32").Value = "  -2.19%  "
$ws.Range("D33").Value = "'4.368"
$ws.Range("E33").Value = "  -2.20%  "
$ws.Range("D34").Value = "'0.04864"
$ws.Range("E34").Value = "  -2.75%  "
$ws.Range("D35").Value = "'0.7580"
$ws.Range("E35").Value = "  -0.32%  "
$ws.Range("D36").Value = "'1.156"
$ws.Range("E36").Value = "  -0.70%  "
$ws.Range("D37").Value = "'0.9997"
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").Value = "'2.724"
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("D39").Value = "'0.01997"
$ws.Range("E39").Value = "  -2.10%  "
$ws.Range("D40").Value = "'2.663"
$ws.Range("E40").Value = "  -1.51%  "
$ws.Range("D41").Value = "'6.496"
$ws.Range("E41").Value = "  +1.01%  "
$ws.Range("D42").Value = "'77.36"
$ws.Range("E42").Value = "  +6.82%  "
$ws.Range("D43").Value = "'2.089"
$ws.Range("E43").Value = "  -2.76%  "
$ws.Range("D44").Value = "'0.9030"
$ws.Range("E44").Value = "  +2.75%  "
$ws.Range("D45").Value = "'0.4415"
$ws.Range("E45").Value = "  -1.54%  "
$ws.Range("D46").Value = "'107.83"
$ws.Range("E46").Value = "  -1.51%  "
$ws.Range("D47").Value = "'0.9984"
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("D48").Value = "'7.776"
$ws.Range("E48").Value = "  +4.05%  "
$ws.Range("D49").Value = "'993.20"
$ws.Range("E49").Value = "  +1.48%  "
$ws.Range("D50").Value = "'0.1249"
$ws.Range("E50").Value = "  -2.52%  "
$ws.Range("D51").Value = "'9.303"
$ws.Range("E51").Value = "  -0.21%  "
